$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "K" (Strike#) values for rows 2-16 (column G)
$newK = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 4
    6  = 4
    7  = 0
    8  = 2
    9  = 6
    10 = 3
    11 = 6
    12 = 3
    13 = 7
    14 = 0
    15 = 3
    16 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
